# Fills in previously-empty "alt data fra Dorte + Eva" measurements on the
# "Analyse_alt" sheet. The dependent %-reduction formulas (row 5, 8, 11, 35,
# 53, ...) already exist in the workbook and simply recalculate once their
# inputs are populated.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Analyse_alt")

# --- Row pair 1 (prøve 1), E=9.25 ---
$ws.Range("F3").Value = 7.84
$ws.Range("G3").Value = 606
$ws.Range("H3").Value = 161
$ws.Range("I3").Value = 1296

$ws.Range("F4").Value = 0.44
$ws.Range("G4").Value = 158
$ws.Range("H4").Value = 198
$ws.Range("I4").Value = 32.5

# --- Row pair 2 (prøve 2), E=9.25 ---
$ws.Range("F6").Value = 11.57
$ws.Range("G6").Value = 360
$ws.Range("H6").Value = 176
$ws.Range("I6").Value = 448
$ws.Range("K6").Value = 465

$ws.Range("F7").Value = 0.85
$ws.Range("G7").Value = 135
$ws.Range("H7").Value = 164
$ws.Range("I7").Value = 27.2
$ws.Range("K7").Value = 180

# --- Row pair 3 (prøve 3), E=9.25 ---
$ws.Range("F9").Value = 8.1
$ws.Range("G9").Value = 235
$ws.Range("H9").Value = 164
$ws.Range("I9").Value = 149
$ws.Range("K9").Value = 455

$ws.Range("F10").Value = 0.91
$ws.Range("G10").Value = 125
$ws.Range("H10").Value = 142
$ws.Range("I10").Value = 5
$ws.Range("K10").Value = 155

# --- Row pair 11 (prøve 11), E=9.5 ---
$ws.Range("F33").Value = 3.89
$ws.Range("G33").Value = 363
$ws.Range("H33").Value = 178
$ws.Range("I33").Value = 463
$ws.Range("K33").Value = 462

$ws.Range("F34").Value = 0.38
$ws.Range("G34").Value = 150
$ws.Range("H34").Value = 178
$ws.Range("I34").Value = 5
$ws.Range("K34").Value = 490

# --- Row pair 17 (prøve 17), E=9.75 ---
$ws.Range("F51").Value = 2.84
$ws.Range("G51").Value = 393
$ws.Range("H51").Value = 177
$ws.Range("I51").Value = 454

$ws.Range("F52").Value = 0.42
$ws.Range("G52").Value = 148
$ws.Range("H52").Value = 158
$ws.Range("I52").Value = 5
$ws.Range("K52").Value = 130

# Put the selection where the author left it when saving.
$ws.Range("M54").Select()
